# Update faturamento_diario for BIBI: correct two May totals and add 4
# more days of May (15-18) that were missing, pushing everything below
# down by 4 rows (Apr/Mar/Fev blocks keep their data, just shifted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the totals already recorded for May 13th and May 14th.
$ws.Range("B14").Value = 26531.86
$ws.Range("B15").Value = 36574.18

# Make room for 4 new rows (May 15-18) right after the current May 14
# row (row 15), pushing the existing April/March/February rows down.
$ws.Rows("16:19").Insert()

# Fill in the newly inserted rows with the May 15-18 data.
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 33940.79
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 2025
$ws.Range("E16").Value = "05/2025"

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 30403.76
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 2025
$ws.Range("E17").Value = "05/2025"

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 14533.8
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 2025
$ws.Range("E18").Value = "05/2025"

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 8085.01
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 2025
$ws.Range("E19").Value = "05/2025"
